# Update the "取得日時" (retrieved at) timestamp for all data rows on the
# "ランサーズ" sheet from 2025-11-30 12:33:44 to 2025-11-30 12:43:19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-11-30 12:33:44"
$newValue = "2025-11-30 12:43:19"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
